$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-20 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-21 Sunday", 2)

$d.Content.Find.Execute("426÷6=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "762÷7=108, 6", 2)
$d.Content.Find.Execute("601÷4=150, 1", $true, $false, $false, $false, $false, $true, 1, $false, "143÷6=23, 5", 2)
$d.Content.Find.Execute("649÷5=129, 4", $true, $false, $false, $false, $false, $true, 1, $false, "145÷7=20, 5", 2)
$d.Content.Find.Execute("401÷7=57, 2", $true, $false, $false, $false, $false, $true, 1, $false, "479÷5=95, 4", 2)
$d.Content.Find.Execute("398÷3=132, 2", $true, $false, $false, $false, $false, $true, 1, $false, "616÷8=77, 0", 2)

$d.Content.Find.Execute("617÷2=308, 1", $true, $false, $false, $false, $false, $true, 1, $false, "589÷3=196, 1", 2)
$d.Content.Find.Execute("895÷3=298, 1", $true, $false, $false, $false, $false, $true, 1, $false, "992÷6=165, 2", 2)
$d.Content.Find.Execute("657÷6=109, 3", $true, $false, $false, $false, $false, $true, 1, $false, "941÷7=134, 3", 2)
$d.Content.Find.Execute("653÷9=72, 5", $true, $false, $false, $false, $false, $true, 1, $false, "812÷7=116, 0", 2)
$d.Content.Find.Execute("571÷4=142, 3", $true, $false, $false, $false, $false, $true, 1, $false, "283÷5=56, 3", 2)

$d.Content.Find.Execute("682÷2=341, 0", $true, $false, $false, $false, $false, $true, 1, $false, "561÷5=112, 1", 2)
$d.Content.Find.Execute("561÷2=280, 1", $true, $false, $false, $false, $false, $true, 1, $false, "456÷3=152, 0", 2)
$d.Content.Find.Execute("600÷8=75, 0", $true, $false, $false, $false, $false, $true, 1, $false, "857÷6=142, 5", 2)
$d.Content.Find.Execute("867÷7=123, 6", $true, $false, $false, $false, $false, $true, 1, $false, "176÷3=58, 2", 2)
$d.Content.Find.Execute("299÷5=59, 4", $true, $false, $false, $false, $false, $true, 1, $false, "455÷3=151, 2", 2)

$d.Content.Find.Execute("877÷5=175, 2", $true, $false, $false, $false, $false, $true, 1, $false, "582÷8=72, 6", 2)
$d.Content.Find.Execute("953÷9=105, 8", $true, $false, $false, $false, $false, $true, 1, $false, "477÷4=119, 1", 2)
$d.Content.Find.Execute("698÷7=99, 5", $true, $false, $false, $false, $false, $true, 1, $false, "754÷5=150, 4", 2)
$d.Content.Find.Execute("866÷6=144, 2", $true, $false, $false, $false, $false, $true, 1, $false, "947÷9=105, 2", 2)
$d.Content.Find.Execute("289÷4=72, 1", $true, $false, $false, $false, $false, $true, 1, $false, "177÷8=22, 1", 2)

$d.Content.Find.Execute("439÷2=219, 1", $true, $false, $false, $false, $false, $true, 1, $false, "679÷2=339, 1", 2)
$d.Content.Find.Execute("819÷5=163, 4", $true, $false, $false, $false, $false, $true, 1, $false, "157÷5=31, 2", 2)
$d.Content.Find.Execute("547÷4=136, 3", $true, $false, $false, $false, $false, $true, 1, $false, "138÷9=15, 3", 2)
$d.Content.Find.Execute("802÷5=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "551÷2=275, 1", 2)
$d.Content.Find.Execute("538÷5=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "792÷3=264, 0", 2)
